$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 4, 6, 7, 8)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "-"
    $ws.Range("E$r").Value = "-"
    $ws.Range("F$r").Value = "MEC-1B-Des. Tec. Mec."
}
